# Atualizando o arquivo XLSX
# Updates a set of betting-odds cell values on Sheet1, row by row,
# matching the committed diff for Jogos_da_Semana_FlashScore_2024-11-05.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 7).Value = 2.35    # G2: 2.2 -> 2.35
$ws.Cells.Item(2, 8).Value = 2.75    # H2: 2.8 -> 2.75
$ws.Cells.Item(2, 9).Value = 3.8     # I2: 4.1 -> 3.8
$ws.Cells.Item(2, 10).Value = 3.25   # J2: 3.2 -> 3.25
$ws.Cells.Item(2, 12).Value = 4.75   # L2: 5 -> 4.75
$ws.Cells.Item(2, 13).Value = 1.18   # M2: 1.17 -> 1.18
$ws.Cells.Item(2, 14).Value = 4.5    # N2: 5 -> 4.5
$ws.Cells.Item(2, 23).Value = 5      # W2: 4.75 -> 5
$ws.Cells.Item(2, 24).Value = 9      # X2: 8.5 -> 9
$ws.Cells.Item(2, 25).Value = 12     # Y2: 11 -> 12
$ws.Cells.Item(2, 26).Value = 23     # Z2: 21 -> 23
$ws.Cells.Item(2, 27).Value = 29     # AA2: 26 -> 29
$ws.Cells.Item(2, 31).Value = 23     # AE2: 26 -> 23
$ws.Cells.Item(2, 33).Value = 7      # AG2: 7.5 -> 7
$ws.Cells.Item(2, 34).Value = 17     # AH2: 19 -> 17
$ws.Cells.Item(2, 35).Value = 15     # AI2: 17 -> 15
$ws.Cells.Item(2, 36).Value = 41     # AJ2: 51 -> 41
$ws.Cells.Item(2, 40).Value = 4      # AN2: 3.75 -> 4
$ws.Cells.Item(2, 41).Value = 17     # AO2: 15 -> 17
$ws.Cells.Item(2, 44).Value = 126    # AR2: 101 -> 126
$ws.Cells.Item(2, 49).Value = 5      # AW2: 5.5 -> 5
$ws.Cells.Item(2, 50).Value = 26     # AX2: 29 -> 26
$ws.Cells.Item(2, 53).Value = 151    # BA2: 201 -> 151

# Row 4
$ws.Cells.Item(4, 19).Value = 1.47   # S4: 1.5 -> 1.47

# Row 5
$ws.Cells.Item(5, 19).Value = 1.3    # S5: 1.33 -> 1.3

# Row 6
$ws.Cells.Item(6, 19).Value = 1.27   # S6: 1.3 -> 1.27

# Row 7
$ws.Cells.Item(7, 7).Value = 2.05    # G7: 2 -> 2.05
$ws.Cells.Item(7, 9).Value = 4       # I7: 4.1 -> 4
$ws.Cells.Item(7, 13).Value = 1.14   # M7: 1.13 -> 1.14
$ws.Cells.Item(7, 14).Value = 5.5    # N7: 6 -> 5.5
$ws.Cells.Item(7, 19).Value = 1.58   # S7: 1.62 -> 1.58
$ws.Cells.Item(7, 36).Value = 41     # AJ7: 51 -> 41
$ws.Cells.Item(7, 50).Value = 26     # AX7: 29 -> 26

# Row 9
$ws.Cells.Item(9, 15).Value = 1.44   # O9: 1.5 -> 1.44
$ws.Cells.Item(9, 16).Value = 2.63   # P9: 2.5 -> 2.63

# Row 10
$ws.Cells.Item(10, 17).Value = 1.6   # Q10: 1.62 -> 1.6
$ws.Cells.Item(10, 18).Value = 2.3   # R10: 2.25 -> 2.3
$ws.Cells.Item(10, 19).Value = 1.3   # S10: 1.27 -> 1.3
$ws.Cells.Item(10, 31).Value = 15    # AE10: 13 -> 15
$ws.Cells.Item(10, 34).Value = 23    # AH10: 26 -> 23

# Row 11
$ws.Cells.Item(11, 19).Value = 1.33  # S11: 1.3 -> 1.33
